# Update to paper 290 april
# The "University of Vermont" closure entry (row 40 — Connor Gage death)
# is removed from the closure-reasons table; the rows below it shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(40).Delete()
